# Fruta / hortaliza, semanal
# Inserts two new weekly price rows for "Femacal de La Calera - Frutilla"
# at the top of the data block (new rows 57-58), pushing the existing
# rows 57-135 down to 59-137.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 57-58; everything at/after row 57 shifts down by 2
# (old row 57 -> new row 59, ..., old row 135 -> new row 137).
$ws.Range("A57:A58").EntireRow.Insert()

# New row 57 - same dimensions/labels as the surrounding data, new weekly figures.
$ws.Range("A57").Value = 3
$ws.Range("B57").Value = "Femacal de La Calera"
$ws.Range("C57").Value = "Coquimbo"
$ws.Range("D57").Value = 44495
$ws.Range("E57").Value = 5
$ws.Range("F57").Value = "Fruta"
$ws.Range("G57").Value = 100101
$ws.Range("H57").Value = "Berries"
$ws.Range("I57").Value = 100112025
$ws.Range("J57").Value = "Frutilla"
$ws.Range("K57").Value = "Sin especificar"
$ws.Range("L57").Value = "Especial"
$ws.Range("M57").Value = 112
$ws.Range("N57").Value = 7000
$ws.Range("O57").Value = 7000
$ws.Range("P57").Value = 7000
$ws.Range("Q57").Value = "$/bandeja 7 kilos"
$ws.Range("R57").Value = "Provincia de Melipilla"
$ws.Range("S57").Value = 1000
$ws.Range("T57").Value = 7

# New row 58.
$ws.Range("A58").Value = 3
$ws.Range("B58").Value = "Femacal de La Calera"
$ws.Range("C58").Value = "Coquimbo"
$ws.Range("D58").Value = 44495
$ws.Range("E58").Value = 5
$ws.Range("F58").Value = "Fruta"
$ws.Range("G58").Value = 100101
$ws.Range("H58").Value = "Berries"
$ws.Range("I58").Value = 100112025
$ws.Range("J58").Value = "Frutilla"
$ws.Range("K58").Value = "Sin especificar"
$ws.Range("L58").Value = "Segunda"
$ws.Range("M58").Value = 75
$ws.Range("N58").Value = 5000
$ws.Range("O58").Value = 5000
$ws.Range("P58").Value = 5000
$ws.Range("Q58").Value = "$/bandeja 7 kilos"
$ws.Range("R58").Value = "Provincia de Melipilla"
$ws.Range("S58").Value = 714
$ws.Range("T58").Value = 7
